# Applies:
#  1. Bold + 14pt (sz/szCs 28) formatting to the "Tracing Sleuth-Zipkin" paragraph.
#  2. Moves the lone "_GoBack" bookmark from the end of the document to sit
#     right after the "In BOM pom.xml add" run (collapsed bookmark).

$d = $word.ActiveDocument

# --- 1. Bold + size the "Tracing Sleuth-Zipkin" heading paragraph ---------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Tracing Sleuth-Zipkin`r") {
        $pr = $p.Range
        $pr.Font.Bold = $true
        $pr.Font.Size = 14
    }
}

# --- 2. Relocate the "_GoBack" bookmark ------------------------------------
# Find the end of the "In BOM pom.xml add" run. Collapsed ranges that land
# exactly on a paragraph boundary confuse Bookmarks.Add, so we temporarily
# append a marker character, anchor the bookmark just before it, then strip
# the marker back out.
$r = $d.Content
$r.Find.Execute("In BOM pom.xml add", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertPos = $r.End
$r.InsertAfter("@")

$bmRange = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$markerRange = $d.Range($insertPos, $insertPos + 1)
$markerRange.Text = ""
